$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A:B").Insert()
$ws.Range("A1").Value = "央管代碼"
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("A1:B1").Merge()
$ws.Range("B1").Borders(10).Weight = -4138
